$d = $word.ActiveDocument

# 1. Title: "ASSIGNMENT - CLUSTERING" -> "CLUSTERING"
$d.Content.Find.Execute("ASSIGNMENT - CLUSTERING", $false, $false, $false, $false, $false, $true, 1, $false, "CLUSTERING", 2)

# 2. Subtitle: "USE python for doing this Assignment" -> "USE Python for customer-segment-analysis"
$d.Content.Find.Execute("USE python for doing this Assignment", $false, $false, $false, $false, $false, $true, 1, $false, "USE Python for customer-segment-analysis", 2)

# 3. Remove the trailing paragraphs: the _GoBack bookmark paragraph, the blank
#    ListParagraph, the "Please share the code..." paragraph, the
#    "Create a ppt..." paragraph, and the final blank paragraph.
$startPar = $d.Paragraphs(8)
$endPar = $d.Paragraphs(12)
$rng = $d.Range($startPar.Range.Start, $endPar.Range.End)
$rng.Delete()
